# Update "想去人数" (number of people interested) in column F
# for the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1217
$ws1.Range("F4").Value = 1128
$ws1.Range("F5").Value = 2474
$ws1.Range("F6").Value = 8034
$ws1.Range("F7").Value = 953
$ws1.Range("F8").Value = 495
$ws1.Range("F11").Value = 466
$ws1.Range("F13").Value = 185
$ws1.Range("F14").Value = 8417
$ws1.Range("F16").Value = 1471
$ws1.Range("F17").Value = 172
$ws1.Range("F20").Value = 220
$ws1.Range("F22").Value = 220
$ws1.Range("F23").Value = 168
$ws1.Range("F25").Value = 124
$ws1.Range("F27").Value = 448
$ws1.Range("F28").Value = 1196
$ws1.Range("F29").Value = 108
$ws1.Range("F30").Value = 69
$ws1.Range("F31").Value = 119
$ws1.Range("F32").Value = 80
$ws1.Range("F33").Value = 108
$ws1.Range("F34").Value = 57
$ws1.Range("F35").Value = 95
$ws1.Range("F36").Value = 87

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1217
$ws4.Range("F4").Value = 1128
$ws4.Range("F5").Value = 2474
$ws4.Range("F6").Value = 8034
$ws4.Range("F7").Value = 953
$ws4.Range("F8").Value = 496
$ws4.Range("F11").Value = 466
$ws4.Range("F13").Value = 185
$ws4.Range("F14").Value = 8417
$ws4.Range("F16").Value = 1471
$ws4.Range("F17").Value = 172
$ws4.Range("F20").Value = 220
$ws4.Range("F22").Value = 220
$ws4.Range("F23").Value = 168
$ws4.Range("F25").Value = 124
$ws4.Range("F27").Value = 448
$ws4.Range("F28").Value = 1196
$ws4.Range("F29").Value = 108
$ws4.Range("F30").Value = 69
$ws4.Range("F31").Value = 119
$ws4.Range("F32").Value = 80
$ws4.Range("F33").Value = 108
$ws4.Range("F34").Value = 57
$ws4.Range("F35").Value = 95
$ws4.Range("F36").Value = 87
